$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates (rich text runs) ---
# A8: "Volume 32   Number  34" -> "...35"  (chars 21-22 = "34")
$ws.Range("A8").Characters(21, 2).Text = "35"

# C9: "Report Covering the Week  8/18/2025  Through  8/24/2025"
#   chars 27-35 = "8/18/2025" -> "8/25/2025"
#   chars 47-55 = "8/24/2025" -> "8/31/2025"
$ws.Range("C9").Characters(27, 9).Text = "8/25/2025"
$ws.Range("C9").Characters(47, 9).Text = "8/31/2025"

# --- Cells whose number format/style must change from a "no data" placeholder
#     to a real numeric style (or vice versa). We use copy/paste-special so the
#     destination cell picks up the exact same style (and, for text, the same
#     shared-string entry) as an existing donor cell elsewhere on the sheet.

$ws.Range("I14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1

$ws.Range("I14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1

$ws.Range("I14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = 1

$ws.Range("L14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100

$ws.Range("I14").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("G29").Value = 1

$ws.Range("L14").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H29").Value = 0

$ws.Range("I14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1

$ws.Range("L14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100

$ws.Range("I14").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1

$ws.Range("L14").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = 0

$ws.Range("I14").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("C33").Value = 1

# C28: numeric 2 -> text placeholder "0" (style reverts to the General "no data" style)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("N14").Value = -94.117647058823
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 27
$ws.Range("K15").Value = 3.846153846153
$ws.Range("L15").Value = 22.727272727272
$ws.Range("M15").Value = 125
$ws.Range("N15").Value = -18.181818181818
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -30
$ws.Range("I16").Value = 128
$ws.Range("J16").Value = 208
$ws.Range("K16").Value = -38.461538461538
$ws.Range("L16").Value = -28.089887640449
$ws.Range("M16").Value = -39.336492890995
$ws.Range("N16").Value = -89.097103918228
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -61.538461538461
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = -30
$ws.Range("I17").Value = 390
$ws.Range("J17").Value = 396
$ws.Range("K17").Value = -1.515151515151
$ws.Range("L17").Value = 18.90243902439
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -26.966292134831
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 17.647058823529
$ws.Range("I18").Value = 148
$ws.Range("J18").Value = 136
$ws.Range("K18").Value = 8.823529411764
$ws.Range("L18").Value = -3.267973856209
$ws.Range("M18").Value = -54.461538461538
$ws.Range("N18").Value = -89.856065798492
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 70
$ws.Range("G19").Value = 86
$ws.Range("H19").Value = -18.60465116279
$ws.Range("I19").Value = 603
$ws.Range("J19").Value = 665
$ws.Range("K19").Value = -9.323308270676
$ws.Range("L19").Value = -1.951219512195
$ws.Range("M19").Value = 36.117381489842
$ws.Range("N19").Value = -16.36615811373
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -37.5
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 33
$ws.Range("H20").Value = -21.212121212121
$ws.Range("I20").Value = 224
$ws.Range("J20").Value = 244
$ws.Range("K20").Value = -8.196721311475
$ws.Range("L20").Value = -20.567375886524
$ws.Range("M20").Value = 22.404371584699
$ws.Range("N20").Value = -87.479038569033
$ws.Range("C21").Value = 41
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = -14.583333333333
$ws.Range("F21").Value = 161
$ws.Range("G21").Value = 197
$ws.Range("H21").Value = -18.274111675126
$ws.Range("I21").Value = 1521
$ws.Range("J21").Value = 1675
$ws.Range("K21").Value = -9.194029850746
$ws.Range("L21").Value = -3.855878634639
$ws.Range("M21").Value = 10.537790697674
$ws.Range("N21").Value = -73.441592456783
$ws.Range("D22").Value = 2
$ws.Range("J22").Value = 28
$ws.Range("K22").Value = -10.714285714285
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 10
$ws.Range("H23").Value = -16.666666666666
$ws.Range("I23").Value = 143
$ws.Range("J23").Value = 149
$ws.Range("K23").Value = -4.026845637583
$ws.Range("L23").Value = -18.75
$ws.Range("M23").Value = 37.5
$ws.Range("C24").Value = 42
$ws.Range("D24").Value = 59
$ws.Range("E24").Value = -28.813559322033
$ws.Range("F24").Value = 197
$ws.Range("G24").Value = 207
$ws.Range("H24").Value = -4.830917874396
$ws.Range("I24").Value = 1484
$ws.Range("J24").Value = 1492
$ws.Range("K24").Value = -0.53619302949
$ws.Range("L24").Value = -6.016466117796
$ws.Range("M24").Value = 46.785361028684
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 39
$ws.Range("E25").Value = -58.974358974359
$ws.Range("F25").Value = 91
$ws.Range("G25").Value = 113
$ws.Range("H25").Value = -19.469026548672
$ws.Range("I25").Value = 751
$ws.Range("J25").Value = 758
$ws.Range("K25").Value = -0.923482849604
$ws.Range("L25").Value = 2.176870748299
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 21
$ws.Range("E26").Value = -14.285714285714
$ws.Range("F26").Value = 61
$ws.Range("G26").Value = 83
$ws.Range("H26").Value = -26.506024096385
$ws.Range("I26").Value = 576
$ws.Range("J26").Value = 645
$ws.Range("K26").Value = -10.697674418604
$ws.Range("L26").Value = 2.127659574468
$ws.Range("M26").Value = -1.369863013698
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 33
$ws.Range("K27").Value = -8.333333333333
$ws.Range("L27").Value = 0
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = -46.153846153846
$ws.Range("J28").Value = 62
$ws.Range("K28").Value = 4.838709677419
$ws.Range("L28").Value = -15.584415584415
$ws.Range("J29").Value = 6
$ws.Range("K29").Value = 50
$ws.Range("L29").Value = -57.142857142857
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -66.666666666666
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = 80
$ws.Range("L30").Value = -52.631578947368
$ws.Range("M30").Value = -40
$ws.Range("N30").Value = -65.384615384615
$ws.Range("F33").Value = 4
$ws.Range("I33").Value = 7
$ws.Range("K33").Value = 133.333333333333
$ws.Range("L33").Value = 40
